# Auto-generated Excel COM-interop script
# Applies updated market-price values scraped by the scheduled runner
# to the Sheets workbook, matching the authoritative commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1099.1875
$ws.Range("I18").Value = 756.9286
$ws.Range("J18").Value = 3495
$ws.Range("K18").Value = 756.9286
$ws.Range("L18").Value = 3495
$ws.Range("M18").Value = -472.9286
$ws.Range("N18").Value = -4063
$ws.Range("H69").Value = 3620
$ws.Range("J69").Value = 3225.7144
$ws.Range("L69").Value = 9677.143199999999
$ws.Range("N69").Value = -11425.1432
$ws.Range("H70").Value = 1281.9286
$ws.Range("I70").Value = 1480
$ws.Range("K70").Value = 4440
$ws.Range("M70").Value = -4170
$ws.Range("H72").Value = 3620
$ws.Range("J72").Value = 3225.7144
$ws.Range("L72").Value = 29031.4296
$ws.Range("N72").Value = -37767.4296
$ws.Range("H73").Value = 1281.9286
$ws.Range("I73").Value = 1480
$ws.Range("K73").Value = 4440
$ws.Range("M73").Value = -3504
$ws.Range("H86").Value = 36475.137
$ws.Range("I86").Value = 45372.74
$ws.Range("J86").Value = 2367.6667
$ws.Range("K86").Value = 45372.74
$ws.Range("L86").Value = 2367.6667
$ws.Range("M86").Value = -44249.74
$ws.Range("N86").Value = -4613.6667
$ws.Range("H89").Value = 36475.137
$ws.Range("I89").Value = 45372.74
$ws.Range("J89").Value = 2367.6667
$ws.Range("K89").Value = 226863.7
$ws.Range("L89").Value = 11838.3335
$ws.Range("M89").Value = -221247.7
$ws.Range("N89").Value = -23070.3335
$ws.Range("H113").Value = 113309.445
$ws.Range("I113").Value = 144886.42
$ws.Range("J113").Value = 2790
$ws.Range("K113").Value = 144886.42
$ws.Range("L113").Value = 2790
$ws.Range("M113").Value = -141632.42
$ws.Range("N113").Value = -9298
$ws.Range("H135").Value = 1619.2858
$ws.Range("I135").Value = 517.7895
$ws.Range("K135").Value = 4660.1055
$ws.Range("M135").Value = -2125.1055

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 22613.262
$ws.Range("I2").Value = 978.1142599999999
$ws.Range("J2").Value = 91452.37
$ws.Range("K2").Value = 978.1142599999999
$ws.Range("L2").Value = 91452.37
$ws.Range("M2").Value = -865.1142599999999
$ws.Range("N2").Value = -91678.37
$ws.Range("H32").Value = 140884.5
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H45").Value = 71217.47
$ws.Range("I45").Value = 92602.09
$ws.Range("K45").Value = 92602.09
$ws.Range("M45").Value = -92225.09
$ws.Range("H97").Value = 23255.889
$ws.Range("I97").Value = 31891.344
$ws.Range("J97").Value = 1999.3846
$ws.Range("K97").Value = 31891.344
$ws.Range("L97").Value = 1999.3846
$ws.Range("M97").Value = -31395.344
$ws.Range("N97").Value = -2991.3846
$ws.Range("H116").Value = 22613.262
$ws.Range("I116").Value = 978.1142599999999
$ws.Range("J116").Value = 91452.37
$ws.Range("K116").Value = 978.1142599999999
$ws.Range("L116").Value = 91452.37
$ws.Range("M116").Value = 1315.88574
$ws.Range("N116").Value = -96040.37
$ws.Range("H122").Value = 1521.2285
$ws.Range("I122").Value = 1423
$ws.Range("J122").Value = 1735.5454
$ws.Range("K122").Value = 4269
$ws.Range("L122").Value = 5206.6362
$ws.Range("M122").Value = -1819
$ws.Range("N122").Value = -10106.6362

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 22613.262
$ws.Range("I3").Value = 978.1142599999999
$ws.Range("J3").Value = 91452.37
$ws.Range("K3").Value = 978.1142599999999
$ws.Range("L3").Value = 91452.37
$ws.Range("M3").Value = -864.1142599999999
$ws.Range("N3").Value = -91680.37
$ws.Range("H20").Value = 48043.684
$ws.Range("I20").Value = 58276.445
$ws.Range("J20").Value = 1996.25
$ws.Range("K20").Value = 58276.445
$ws.Range("L20").Value = 1996.25
$ws.Range("M20").Value = -58029.445
$ws.Range("N20").Value = -2490.25
$ws.Range("H80").Value = 2810.6667
$ws.Range("I80").Value = 845.7778
$ws.Range("K80").Value = 845.7778
$ws.Range("M80").Value = 152.2222
$ws.Range("H83").Value = 2810.6667
$ws.Range("I83").Value = 845.7778
$ws.Range("K83").Value = 4228.889
$ws.Range("M83").Value = 763.1109999999999
$ws.Range("H94").Value = 753.2143
$ws.Range("I94").Value = 734.2308
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 734.2308
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -283.2308
$ws.Range("N94").Value = -1902

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 77491.664
$ws.Range("J138").Value = 77491.664
$ws.Range("L138").Value = 77491.664
$ws.Range("N138").Value = -87771.664

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 609.0769
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 657.3333
$ws.Range("K23").Value = 90
$ws.Range("L23").Value = 1971.9999
$ws.Range("M23").Value = 145
$ws.Range("N23").Value = -2441.9999
$ws.Range("H98").Value = 100998
$ws.Range("I98").Value = 351.5
$ws.Range("J98").Value = 126159.625
$ws.Range("K98").Value = 1054.5
$ws.Range("L98").Value = 378478.875
$ws.Range("M98").Value = 443.5
$ws.Range("N98").Value = -381474.875
$ws.Range("H131").Value = 1188.3152
$ws.Range("J131").Value = 1217.784
$ws.Range("L131").Value = 3653.352
$ws.Range("N131").Value = -13733.352

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1895.909
$ws.Range("J113").Value = 1895
$ws.Range("L113").Value = 1895
$ws.Range("N113").Value = -6235

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2092.8572
$ws.Range("I22").Value = 3550
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3550
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -3255
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 2092.8572
$ws.Range("I27").Value = 3550
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 3550
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -3443
$ws.Range("N27").Value = -1214
$ws.Range("H55").Value = 695.7560999999999
$ws.Range("I55").Value = 1116.5333
$ws.Range("J55").Value = 453
$ws.Range("K55").Value = 1116.5333
$ws.Range("L55").Value = 453
$ws.Range("M55").Value = -943.5333000000001
$ws.Range("N55").Value = -799
$ws.Range("H68").Value = 3832.7273
$ws.Range("I68").Value = 1596.6666
$ws.Range("J68").Value = 4671.25
$ws.Range("K68").Value = 1596.6666
$ws.Range("L68").Value = 4671.25
$ws.Range("M68").Value = -847.6666
$ws.Range("N68").Value = -6169.25
$ws.Range("H71").Value = 3832.7273
$ws.Range("I71").Value = 1596.6666
$ws.Range("J71").Value = 4671.25
$ws.Range("K71").Value = 7983.333000000001
$ws.Range("L71").Value = 23356.25
$ws.Range("M71").Value = -4239.333000000001
$ws.Range("N71").Value = -30844.25
$ws.Range("H132").Value = 5678.1055
$ws.Range("I132").Value = 6207.5
$ws.Range("J132").Value = 4770.5713
$ws.Range("K132").Value = 18622.5
$ws.Range("L132").Value = 14311.7139
$ws.Range("M132").Value = -16092.5
$ws.Range("N132").Value = -19371.7139

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1930.7941
$ws.Range("I132").Value = 2034.62
$ws.Range("K132").Value = 6103.86
$ws.Range("M132").Value = -3573.86

